$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 8118.727
$ws.Range("I28").Value = 3059.5789
$ws.Range("J28").Value = 40160
$ws.Range("K28").Value = 3059.5789
$ws.Range("L28").Value = 40160
$ws.Range("M28").Value = -2574.5789
$ws.Range("N28").Value = -41130

$ws.Range("H45").Value = 1758.5
$ws.Range("I45").Value = 1017
$ws.Range("K45").Value = 3051
$ws.Range("M45").Value = -2859

$ws.Range("H48").Value = 722.2222
$ws.Range("J48").Value = 722.2222
$ws.Range("L48").Value = 2166.6666
$ws.Range("N48").Value = -2750.6666

$ws.Range("H49").Value = 75
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").ClearContents()

$ws.Range("H56").Value = 722.2222
$ws.Range("J56").Value = 722.2222
$ws.Range("L56").Value = 2166.6666
$ws.Range("N56").Value = -3234.6666

$ws.Range("H62").Value = 1801.3636
$ws.Range("I62").Value = 970
$ws.Range("J62").Value = 2113.125
$ws.Range("K62").Value = 970
$ws.Range("L62").Value = 2113.125
$ws.Range("M62").Value = -346
$ws.Range("N62").Value = -3361.125

$ws.Range("H65").Value = 1801.3636
$ws.Range("I65").Value = 970
$ws.Range("J65").Value = 2113.125
$ws.Range("K65").Value = 4850
$ws.Range("L65").Value = 10565.625
$ws.Range("M65").Value = -1730
$ws.Range("N65").Value = -16805.625

$ws.Range("H116").Value = 3066.6667
$ws.Range("I116").Value = 2480
$ws.Range("K116").Value = 2480
$ws.Range("M116").Value = 962

$ws.Range("H132").Value = 4147.68
$ws.Range("I132").Value = 4034.4783
$ws.Range("K132").Value = 12103.4349
$ws.Range("M132").Value = -9573.4349

$ws.Range("H135").Value = 3536.2222
$ws.Range("I135").Value = 2358
$ws.Range("J135").Value = 5009
$ws.Range("K135").Value = 21222
$ws.Range("L135").Value = 45081
$ws.Range("M135").Value = -18687
$ws.Range("N135").Value = -50151

$ws.Range("H138").Value = 2496.9324
$ws.Range("I138").Value = 1545.8611
$ws.Range("J138").Value = 3397.9473
$ws.Range("K138").Value = 4637.5833
$ws.Range("L138").Value = 10193.8419
$ws.Range("M138").Value = 502.4166999999998
$ws.Range("N138").Value = -20473.8419

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2874.606
$ws.Range("I61").Value = 2317.0833
$ws.Range("J61").Value = 4361.3335
$ws.Range("K61").Value = 2317.0833
$ws.Range("L61").Value = 4361.3335
$ws.Range("M61").Value = -2105.0833
$ws.Range("N61").Value = -4785.3335

$ws.Range("H74").Value = 3050.5386
$ws.Range("I74").Value = 2250
$ws.Range("J74").Value = 3550.875
$ws.Range("K74").Value = 2250
$ws.Range("L74").Value = 3550.875
$ws.Range("M74").Value = -1376
$ws.Range("N74").Value = -5298.875

$ws.Range("H77").Value = 3050.5386
$ws.Range("I77").Value = 2250
$ws.Range("J77").Value = 3550.875
$ws.Range("K77").Value = 11250
$ws.Range("L77").Value = 17754.375
$ws.Range("M77").Value = -6882
$ws.Range("N77").Value = -26490.375

$ws.Range("H122").Value = 335333
$ws.Range("I122").Value = 500499.5
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 1501498.5
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -1499048.5
$ws.Range("N122").Value = -19900

$ws.Range("H136").Value = 2874.606
$ws.Range("I136").Value = 2317.0833
$ws.Range("J136").Value = 4361.3335
$ws.Range("K136").Value = 6951.249899999999
$ws.Range("L136").Value = 13084.0005
$ws.Range("M136").Value = -4401.249899999999
$ws.Range("N136").Value = -18184.0005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H93").Value = 70000
$ws.Range("J93").Value = 70000
$ws.Range("L93").Value = 70000
$ws.Range("N93").Value = -73744

$ws.Range("H102").Value = 20006.125
$ws.Range("I102").Value = 8578.429
$ws.Range("K102").Value = 8578.429
$ws.Range("M102").Value = -5333.429

$ws.Range("H134").Value = 8169
$ws.Range("I134").Value = 10000
$ws.Range("J134").Value = 7253.5
$ws.Range("K134").Value = 30000
$ws.Range("L134").Value = 21760.5
$ws.Range("M134").Value = -27465
$ws.Range("N134").Value = -26830.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6369.4116
$ws.Range("I31").Value = 1178.4
$ws.Range("J31").Value = 11360.77
$ws.Range("K31").Value = 1178.4
$ws.Range("L31").Value = 11360.77
$ws.Range("M31").Value = -883.4000000000001
$ws.Range("N31").Value = -11950.77

$ws.Range("H34").Value = 6369.4116
$ws.Range("I34").Value = 1178.4
$ws.Range("J34").Value = 11360.77
$ws.Range("K34").Value = 1178.4
$ws.Range("L34").Value = 11360.77
$ws.Range("M34").Value = -976.4000000000001
$ws.Range("N34").Value = -11764.77

$ws.Range("H122").Value = 1905.2858
$ws.Range("I122").Value = 1252.75
$ws.Range("J122").Value = 2058.8235
$ws.Range("K122").Value = 3758.25
$ws.Range("L122").Value = 6176.470499999999
$ws.Range("M122").Value = -1308.25
$ws.Range("N122").Value = -11076.4705

$ws.Range("H132").Value = 8335783.5
$ws.Range("I132").Value = 1863.3572
$ws.Range("J132").Value = 27781598
$ws.Range("K132").Value = 5590.071599999999
$ws.Range("L132").Value = 83344794
$ws.Range("M132").Value = -3060.071599999999
$ws.Range("N132").Value = -83349854

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 892.41174
$ws.Range("J5").Value = 1328.7142
$ws.Range("L5").Value = 3986.1426
$ws.Range("N5").Value = -4210.142599999999

$ws.Range("H36").Value = 4166.6665
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 4166.6665
$ws.Range("K36").Value = 0
$ws.Range("L36").Value = 12499.9995
$ws.Range("M36").ClearContents()
$ws.Range("N36").Value = -12837.9995

$ws.Range("H50").Value = 83333840
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 83333840
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 250001520
$ws.Range("M50").ClearContents()
$ws.Range("N50").Value = -250002482

$ws.Range("H53").Value = 83333840
$ws.Range("I53").Value = 0
$ws.Range("J53").Value = 83333840
$ws.Range("K53").Value = 0
$ws.Range("L53").Value = 250001520
$ws.Range("M53").ClearContents()
$ws.Range("N53").Value = -250002482

$ws.Range("H98").Value = 600.5
$ws.Range("I98").Value = 701
$ws.Range("K98").Value = 2103
$ws.Range("M98").Value = -605

$ws.Range("H122").Value = 5363.5454
$ws.Range("I122").Value = 500
$ws.Range("J122").Value = 9416.5
$ws.Range("K122").Value = 4500
$ws.Range("L122").Value = 84748.5
$ws.Range("M122").Value = -2050
$ws.Range("N122").Value = -89648.5

$ws.Range("H135").Value = 892.41174
$ws.Range("J135").Value = 1328.7142
$ws.Range("L135").Value = 11958.4278
$ws.Range("N135").Value = -17028.4278

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 367.66666
$ws.Range("I2").Value = 53.666668
$ws.Range("J2").Value = 681.6667
$ws.Range("K2").Value = 53.666668
$ws.Range("L2").Value = 681.6667
$ws.Range("M2").Value = 59.333332
$ws.Range("N2").Value = -907.6667

$ws.Range("H102").Value = 1583.2632
$ws.Range("I102").Value = 1447.7142
$ws.Range("J102").Value = 1962.8
$ws.Range("K102").Value = 1447.7142
$ws.Range("L102").Value = 1962.8
$ws.Range("M102").Value = 174.2858000000001
$ws.Range("N102").Value = -5206.8

$ws.Range("H122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 8886.286
$ws.Range("I22").Value = 2090.818
$ws.Range("J22").Value = 13283.353
$ws.Range("K22").Value = 2090.818
$ws.Range("L22").Value = 13283.353
$ws.Range("M22").Value = -1795.818
$ws.Range("N22").Value = -13873.353

$ws.Range("H27").Value = 8886.286
$ws.Range("I27").Value = 2090.818
$ws.Range("J27").Value = 13283.353
$ws.Range("K27").Value = 2090.818
$ws.Range("L27").Value = 13283.353
$ws.Range("M27").Value = -1983.818
$ws.Range("N27").Value = -13497.353

$ws.Range("H40").Value = 336334.66
$ws.Range("I40").Value = 502502
$ws.Range("K40").Value = 502502
$ws.Range("M40").Value = -502366

$ws.Range("H46").Value = 2300
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 2300
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 2300
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -2676

$ws.Range("H61").Value = 4230.8
$ws.Range("I61").Value = 4452
$ws.Range("K61").Value = 4452
$ws.Range("M61").Value = -4250

$ws.Range("H68").Value = 1540.5883
$ws.Range("I68").Value = 1483.8462
$ws.Range("K68").Value = 1483.8462
$ws.Range("M68").Value = -734.8462

$ws.Range("H71").Value = 1540.5883
$ws.Range("I71").Value = 1483.8462
$ws.Range("K71").Value = 7419.231
$ws.Range("M71").Value = -3675.231

$ws.Range("H96").Value = 85000
$ws.Range("J96").Value = 85000
$ws.Range("L96").Value = 85000
$ws.Range("N96").Value = -90492

$ws.Range("H100").Value = 3940.75
$ws.Range("I100").Value = 3921
$ws.Range("J100").Value = 4000
$ws.Range("K100").Value = 3921
$ws.Range("L100").Value = 4000
$ws.Range("M100").Value = -3380
$ws.Range("N100").Value = -5082

$ws.Range("H113").Value = 4230.8
$ws.Range("I113").Value = 4452
$ws.Range("K113").Value = 4452
$ws.Range("M113").Value = -2282

$ws.Range("H122").Value = 3500
$ws.Range("I122").Value = 2000
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 6000
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -3550
$ws.Range("N122").Value = -19900

$ws.Range("H132").Value = 2886.2222
$ws.Range("I132").Value = 1998.4762
$ws.Range("K132").Value = 5995.4286
$ws.Range("M132").Value = -3465.4286

$ws.Range("H136").Value = 16669160
$ws.Range("J136").Value = 23811428
$ws.Range("L136").Value = 71434284
$ws.Range("N136").Value = -71439384

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 726.1111
$ws.Range("I107").Value = 729.375
$ws.Range("K107").Value = 2188.125
$ws.Range("M107").Value = -268.125

$ws.Range("H136").Value = 4281.8
$ws.Range("I136").Value = 3580.8462
$ws.Range("J136").Value = 5583.5713
$ws.Range("K136").Value = 10742.5386
$ws.Range("L136").Value = 16750.7139
$ws.Range("M136").Value = -8192.5386
$ws.Range("N136").Value = -21850.7139
